# Apply the edits described by the diff to "Hoja1" (sheet1.xml)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: date serial changed from 45406 (2024-04-24) to 45436 (2024-05-24)
$ws.Range("A1").Value = 45436

# D29: price changed from 410 to 1230
$ws.Range("D29").Value = 1230

# D30: price changed from 445 to 1290
$ws.Range("D30").Value = 1290
